# Grievance Policy and Procedure.docx -- update on 03/05/2021 at 17:24
#
# 1) "Written: Spring Term 2020"            -> "Written: Summer Term 2021"
# 2) "Date of Next review: Spring Term 2021" -> "Date of Next review: summer Term 2022"
# 3) "12th April 2020"                       -> "15th April 2021"   (keep "th" superscript run untouched)
# 4) "This policy aims to enable employees to raise concerns ..."
#      -> "... employees, contractors and supply staff (herein referred to as employees), to raise concerns ..."
# 5) "... according to GDPR2016/679 and the school's Data Management Policy"
#      -> "... according to UK General Data Protection Regulation (UK GDPR), tailored by the
#          Data Protection Act 2018, and the school's Data Management Policy"
# 6) Remove the stray "_GoBack" bookmark left around "...within|5 working days..."

$d = $word.ActiveDocument

# 1) Written: Spring Term 2020 -> Written: Summer Term 2021
$d.Content.Find.Execute("Spring Term 2020", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Summer Term 2021", 2) | Out-Null

# 2) Date of Next review: Spring Term 2021 -> Date of Next review: summer Term 2022
$d.Content.Find.Execute("of Next review: Spring Term 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "of Next review: summer Term 2022", 2) | Out-Null

# 3) 12th April 2020 -> 15th April 2021 (done as two pieces so the superscript "th" run is untouched)
#    "12" and "April 2020" are each unique in the document, so plain (non-whole-word) matches are safe.
$d.Content.Find.Execute("12", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "15", 2) | Out-Null
$d.Content.Find.Execute("April 2020", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "April 2021", 2) | Out-Null

# 4) Add ", contractors and supply staff (herein referred to as employees)," after "employees"
$d.Content.Find.Execute("enable employees to raise concerns", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "enable employees, contractors and supply staff (herein referred to as employees), to raise concerns", `
                         2) | Out-Null

# 5) Replace the old "GDPR2016/679" reference with the UK GDPR / Data Protection Act 2018 wording
$d.Content.Find.Execute("according to GDPR2016/679 and the school", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "according to UK General Data Protection Regulation (UK GDPR), tailored by the Data Protection Act 2018, and the school", `
                         2) | Out-Null

# 6) Drop the leftover "_GoBack" bookmark (touch the text spanning across it, without changing the text)
$d.Content.Find.Execute("Panel) within 5", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Panel) within 5", 2) | Out-Null
